# menambahkan & mengurutkan nota berdasarkan jenis barang
# Adds a new "Jenis Barang" (item type) header column so the import
# template can be sorted/grouped by item type.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The header row (row 4) currently runs from A4 to K4; add the new
# "Jenis Barang" header in the next column, L4.
$ws.Cells.Item(4, 12).Value = "Jenis Barang"

# Give the new header cell the same look (font/fill) as the rest of the
# header row before centering everything.
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Center every header cell horizontally (A4:L4), matching the rest of
# the sheet's "Jenis Barang" addition.
$ws.Range("A4:L4").HorizontalAlignment = -4108   # xlCenter

# Column L should match the width used by the other data columns (H:I).
$ws.Columns.Item(12).ColumnWidth = $ws.Columns.Item(9).ColumnWidth

# Reflect the new active cell now that the header row was extended.
$ws.Range("L5").Select()
